$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.228.47"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.908.33"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.97%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.26"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5246"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3785"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.31%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07278"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.28"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8993"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07686"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.913.30"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.09"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.276"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008655"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.55"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.296.85"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.090"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.153.52"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.64"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.445"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.330"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +11.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.74"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.18"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.735"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.38%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.87"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.969"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +5.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.814"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.43%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05074"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7965"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.41%  "

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.245"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.986"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.60%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.308"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.611"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5695"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01995"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.076"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.25%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.023"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.11%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.644"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.25"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.84%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1520"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.85%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4865"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.27"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.43%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.612"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.57"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.01"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.33%  "
